$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 2 and 3 (PLK2 <-> THBS1), updating B column values per diff
$ws.Range("A2").Value = "THBS1"
$ws.Range("B2").Value = "Associated with Dysbiosis, Cancer, Neurological Diseases, Rare Diseases, Viral Diseases"
$ws.Range("C2").Value = "Bacillus, Fusobacterium"

$ws.Range("A3").Value = "PLK2"
$ws.Range("B3").Value = "Associated with Dysbiosis, Cancer, Neurological Diseases, Rare Diseases"
$ws.Range("C3").Value = "Bacillus, Fusobacterium"

# Add new row 4
$ws.Range("A4").Value = "CXCL8"
$ws.Range("B4").Value = "Associated with Dysbiosis, Cancer, Neurological Diseases, Pathogenic Bacteria, Rare Diseases, Viral Diseases"
$ws.Range("C4").Value = "Bacillus, Lactobacillus"
